# 21 Aug 2023 update: add Tong (Cyber team) as a playing member.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Tong (row 48, team Cyber) as playing this round.
$ws.Range("D48").Value = 1

# Recalculate so dependent formulas (E, L2, M2, ...) pick up the new total.
$excel.Calculate()

# Leave the cursor where the editor ended up after the edit, scrolled back
# to the top of the sheet.
$ws.Range("I7").Select()
